$d = $word.ActiveDocument

# --- Part 1: "RF03 - BUSCAR FUNCIONARIO" paragraph ---------------------
# Original: "...nome ou CPF..."  ->  "...nome parcial ou CPF..."
# The new text "parcial" must land in its own run (matching the source
# edit, which split the single run into three runs with identical
# character formatting).
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("nome ou CPF", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find the funcionarios search-criteria sentence"
}

$insertAt = $rng1.Start + 5   # just after "nome " (before "ou CPF")
$insPoint = $d.Range($insertAt, $insertAt)
$insPoint.InsertBefore("parcial ")

# "parcial" now spans insertAt .. insertAt+7 ; force Word to carve this
# span into its own run by touching a character-formatting property
# (round-tripped back to its original value so the visible formatting
# is unchanged, but the run boundary is created).
$newWord = $d.Range($insertAt, $insertAt + 7)
$newWord.Font.Size = 24
$newWord.Font.Size = 12

# --- Part 2: "RF20 - BUSCAR CLIENTE" paragraph --------------------------
# Originally split across six runs ("...nome " / "parcial" / " ou telefone "
# / "parcial" / ". Caso ... retornar " / "id," / " CPF, nome..."); the
# commit collapses them back into a single run with the same text.
$fullClientText = "O sistema deve permitir buscar clientes com os critérios: CPF, nome parcial ou telefone parcial. Caso nenhum critério seja informado, a busca deve retornar todos os clientes. A busca deve retornar id, CPF, nome, telefone e data de nascimento."

$rng2 = $d.Content
$found2 = $rng2.Find.Execute($fullClientText, $true, $false, $false, $false, $false, $true, 1, $false, $fullClientText, 2)
if (-not $found2) {
    throw "Could not find/replace the clientes search-criteria sentence"
}
